$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as literal text, even when the text looks like
# a number, without leaving a permanent style change on the cell (restores
# the cell's original Style right after assigning the value).
function Set-TextCell($row, $col, [string]$text) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

# --- Simple price (column D) updates ---
Set-TextCell 2  4 "236.52"
Set-TextCell 5  4 "0.05568"
Set-TextCell 6  4 "3.369"
Set-TextCell 7  4 "6.456"
Set-TextCell 8  4 "0.7988"
Set-TextCell 9  4 "1.038"
Set-TextCell 10 4 "0.1398"
Set-TextCell 11 4 "0.07305"
Set-TextCell 12 4 "0.03214"
Set-TextCell 13 4 "0.02912"
Set-TextCell 14 4 "0.09245"
Set-TextCell 15 4 "0.001669"
Set-TextCell 16 4 "3.254"
Set-TextCell 17 4 "0.04758"
Set-TextCell 18 4 "0.0005710"
Set-TextCell 19 4 "0.006261"
Set-TextCell 20 4 "0.005062"
Set-TextCell 21 4 "0.001048"
Set-TextCell 22 4 "0.0001499"
Set-TextCell 23 4 "0.0004182"
Set-TextCell 24 4 "3.956"
Set-TextCell 27 4 "0.1295"
Set-TextCell 40 4 "0.04126"
Set-TextCell 41 4 "0.006964"

# --- Row 42/43: CEJI and BKEXToken swap places (with updated prices) ---
Set-TextCell 42 2 "BKEXToken"
Set-TextCell 42 3 "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell 42 4 "0.1038"
Set-TextCell 42 5 "41BKEXTokenBKK"

Set-TextCell 43 2 "CEJI"
Set-TextCell 43 3 "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextCell 43 4 "0.002919"
Set-TextCell 43 5 "42CEJICEJI"

# --- Remaining simple price (column D) updates ---
Set-TextCell 44 4 "0.008754"
Set-TextCell 45 4 "0.00005441"
Set-TextCell 47 4 "0.6800"
Set-TextCell 48 4 "0.03211"

Write-Host "Applied cryptos.xlsx symbol list update"
